$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.330.17'
$ws.Range("E2").Value = '  -0.73%  '
$ws.Range("D3").Value = '1.827.25'
$ws.Range("E3").Value = '  +1.87%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '329.76'
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4463'
$ws.Range("E7").Value = '  +1.47%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3785'
$ws.Range("E8").Value = '  +1.63%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.79'
$ws.Range("E9").Value = '  -1.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07779'
$ws.Range("E10").Value = '  +2.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.140'
$ws.Range("E11").Value = '  +0.92%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.35'
$ws.Range("E12").Value = '  -0.98%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.003'
$ws.Range("E13").Value = '  -0.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.368'
$ws.Range("E14").Value = '  +2.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.539'
$ws.Range("E15").Value = '  +1.01%  '
$ws.Range("D16").Value = '1.838.77'
$ws.Range("E16").Value = '  +2.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.21'
$ws.Range("E17").Value = '  +16.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001087'
$ws.Range("E18").Value = '  +0.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06402'
$ws.Range("E19").Value = '  -4.29%  '
$ws.Range("E20").Value = '  -0.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.57'
$ws.Range("E21").Value = '  +0.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.370'
$ws.Range("E22").Value = '  +2.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.5419'
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("D24").Value = '28.396.20'
$ws.Range("E24").Value = '  -0.54%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.74'
$ws.Range("E25").Value = '  +0.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.273'
$ws.Range("E26").Value = '  -6.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.85'
$ws.Range("E27").Value = '  +2.43%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '154.81'
$ws.Range("E28").Value = '  +1.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.378'
$ws.Range("E29").Value = '  +2.13%  '
$ws.Range("D30").Value = '2.045.25'
$ws.Range("E30").Value = '  +2.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '129.18'
$ws.Range("E31").Value = '  -0.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.213'
$ws.Range("E32").Value = '  -7.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.940'
$ws.Range("E33").Value = '  +3.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.09322'
$ws.Range("E34").Value = '  +0.82%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.673'
$ws.Range("E35").Value = '  -7.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '13.10'
$ws.Range("E36").Value = '  +8.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02358'
$ws.Range("E37").Value = '  +2.52%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2203'
$ws.Range("E38").Value = '  -1.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6659'
$ws.Range("E39").Value = '  +1.73%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06287'
$ws.Range("E40").Value = '  +0.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.219'
$ws.Range("E41").Value = '  +0.84%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.202'
$ws.Range("E42").Value = '  +0.59%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.171'
$ws.Range("E43").Value = '  +2.61%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.16'
$ws.Range("E44").Value = '  +1.93%  '
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  -0.15%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.400'
$ws.Range("E46").Value = '  -2.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6161'
$ws.Range("E47").Value = '  +1.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.781'
$ws.Range("E48").Value = '  -0.48%  '
$ws.Range("E49").Value = '  +3.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '127.59'
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07019'
$ws.Range("E51").Value = '  +0.35%  '
